$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474993058530/1743637847690/PaZ2GKc6aN-camera-video-segment-driver-1743637850190.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSCQSKQN3K%2F20250403%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250403T210046Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIz%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIA0f1nL%2FXh13bkMdIT2AmwaiBdDQqQS9e0S%2FljvQ4gKaAiEA66C2PeoZXPhcqkfZj%2FSmhFW0JZF358cY4pis%2FqXODcgq5gMI9f%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDFOSUKyEHXsANiPAbiq6A14JdRT1c3F1axdrWHYBaP4dVElP6xpcwzyENXLCxholSW0EyRm2KoLDT8a4n9B9ZJuKgjPEt5LjThVMWWKOg40pb0dlf0FjqleeWpY7EqFLWNUnEERekltQN1FrsnhNaybTLBWDUbK7QqIe9ribi4fvhm8GEKsmOeqS1im8xQ8gTbqHK%2BARSuIYRfmjOrPZwuNh9HBM5HM6Sxc9JNYGsc3xbc%2F59xOniBVMkmvIee3HSFyuAGIGoSi4x4FFBPKhucNC%2BHrhhGHM%2F4oXd8tJ2CxfVpjlsGuV3fwyfd452aPIyuOjvS9Ws7cW7VArjg4xRY%2BVM%2BI0%2Fygju8d9lKYak34Yq%2B2TunJCxHrLWr9XmZ1RMg3SsCOdudtodwvZl3mew2CXGL%2B4MOPnu%2B3u1gZD4aKcpbo3ODmhav8qEerJJapJ4ZfbHF%2FH9zwKEGG3Q7ntq1SjpOOItgE%2B6eQ4byWj1pIkxEEUltRHx8Cie8fmzXi0PYYHPkLcAbOGoJPSu6tnhHtgqv7g0yy%2BSo5mqe7G8wZzy3edkd4BpevZaBsLiyQME6wkh%2BxodKB7wbF%2FhoZJ0VY2CJjX1U6i0qYwy8O7vwY6pQH27dqHdKwnvo1S%2B%2BTEkT6eVRCthiib3yqPaq7st7YzJ0LAAIq6zu7k0iBrQLWobAv1vd337f4fGLbgEpafIwWBZouQ11G%2FKzHkRAZN%2F0AhKRdBlD5gDHSwKcdJdJ5aSyQykz0sYCO%2Bq%2BkMt1upPfVIzb48jgPBVMfe7Vbso6c%2FIrl5F9HJojfAmmTbKgieM40Lq0vvS9Shm8v%2FZIOYPJniMHJakrk%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2004%20Apr%202025%2005%3A00%3A46%20GMT&X-Amz-Signature=a89d2351c6bc9280823a3aea6d4b534db14bfd69cc3d5ca8f357294a6167d81f'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474993058530/1743634721601/5dyCDJLTC3-camera-video-segment-driver-1743634724101.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSCQSKQN3K%2F20250403%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250403T210046Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIz%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIA0f1nL%2FXh13bkMdIT2AmwaiBdDQqQS9e0S%2FljvQ4gKaAiEA66C2PeoZXPhcqkfZj%2FSmhFW0JZF358cY4pis%2FqXODcgq5gMI9f%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDFOSUKyEHXsANiPAbiq6A14JdRT1c3F1axdrWHYBaP4dVElP6xpcwzyENXLCxholSW0EyRm2KoLDT8a4n9B9ZJuKgjPEt5LjThVMWWKOg40pb0dlf0FjqleeWpY7EqFLWNUnEERekltQN1FrsnhNaybTLBWDUbK7QqIe9ribi4fvhm8GEKsmOeqS1im8xQ8gTbqHK%2BARSuIYRfmjOrPZwuNh9HBM5HM6Sxc9JNYGsc3xbc%2F59xOniBVMkmvIee3HSFyuAGIGoSi4x4FFBPKhucNC%2BHrhhGHM%2F4oXd8tJ2CxfVpjlsGuV3fwyfd452aPIyuOjvS9Ws7cW7VArjg4xRY%2BVM%2BI0%2Fygju8d9lKYak34Yq%2B2TunJCxHrLWr9XmZ1RMg3SsCOdudtodwvZl3mew2CXGL%2B4MOPnu%2B3u1gZD4aKcpbo3ODmhav8qEerJJapJ4ZfbHF%2FH9zwKEGG3Q7ntq1SjpOOItgE%2B6eQ4byWj1pIkxEEUltRHx8Cie8fmzXi0PYYHPkLcAbOGoJPSu6tnhHtgqv7g0yy%2BSo5mqe7G8wZzy3edkd4BpevZaBsLiyQME6wkh%2BxodKB7wbF%2FhoZJ0VY2CJjX1U6i0qYwy8O7vwY6pQH27dqHdKwnvo1S%2B%2BTEkT6eVRCthiib3yqPaq7st7YzJ0LAAIq6zu7k0iBrQLWobAv1vd337f4fGLbgEpafIwWBZouQ11G%2FKzHkRAZN%2F0AhKRdBlD5gDHSwKcdJdJ5aSyQykz0sYCO%2Bq%2BkMt1upPfVIzb48jgPBVMfe7Vbso6c%2FIrl5F9HJojfAmmTbKgieM40Lq0vvS9Shm8v%2FZIOYPJniMHJakrk%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2004%20Apr%202025%2005%3A00%3A46%20GMT&X-Amz-Signature=9802e9feb1ca34e2be97365dd61933af005a4d3a2011c1bd2aa3c455e03b48d4'
$ws.Range("K5").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991109864/1743618603072/lsfY7Elc8s-camera-video-segment-driver-1743618608072.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSCQSKQN3K%2F20250403%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250403T210046Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIz%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIA0f1nL%2FXh13bkMdIT2AmwaiBdDQqQS9e0S%2FljvQ4gKaAiEA66C2PeoZXPhcqkfZj%2FSmhFW0JZF358cY4pis%2FqXODcgq5gMI9f%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDFOSUKyEHXsANiPAbiq6A14JdRT1c3F1axdrWHYBaP4dVElP6xpcwzyENXLCxholSW0EyRm2KoLDT8a4n9B9ZJuKgjPEt5LjThVMWWKOg40pb0dlf0FjqleeWpY7EqFLWNUnEERekltQN1FrsnhNaybTLBWDUbK7QqIe9ribi4fvhm8GEKsmOeqS1im8xQ8gTbqHK%2BARSuIYRfmjOrPZwuNh9HBM5HM6Sxc9JNYGsc3xbc%2F59xOniBVMkmvIee3HSFyuAGIGoSi4x4FFBPKhucNC%2BHrhhGHM%2F4oXd8tJ2CxfVpjlsGuV3fwyfd452aPIyuOjvS9Ws7cW7VArjg4xRY%2BVM%2BI0%2Fygju8d9lKYak34Yq%2B2TunJCxHrLWr9XmZ1RMg3SsCOdudtodwvZl3mew2CXGL%2B4MOPnu%2B3u1gZD4aKcpbo3ODmhav8qEerJJapJ4ZfbHF%2FH9zwKEGG3Q7ntq1SjpOOItgE%2B6eQ4byWj1pIkxEEUltRHx8Cie8fmzXi0PYYHPkLcAbOGoJPSu6tnhHtgqv7g0yy%2BSo5mqe7G8wZzy3edkd4BpevZaBsLiyQME6wkh%2BxodKB7wbF%2FhoZJ0VY2CJjX1U6i0qYwy8O7vwY6pQH27dqHdKwnvo1S%2B%2BTEkT6eVRCthiib3yqPaq7st7YzJ0LAAIq6zu7k0iBrQLWobAv1vd337f4fGLbgEpafIwWBZouQ11G%2FKzHkRAZN%2F0AhKRdBlD5gDHSwKcdJdJ5aSyQykz0sYCO%2Bq%2BkMt1upPfVIzb48jgPBVMfe7Vbso6c%2FIrl5F9HJojfAmmTbKgieM40Lq0vvS9Shm8v%2FZIOYPJniMHJakrk%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2004%20Apr%202025%2005%3A00%3A46%20GMT&X-Amz-Signature=2207a1c2115365f932a338388e69fa126c72acbb37e6577ae9b54b1424c492e4'
$ws.Range("L5").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991109864/1743618603072/krpw3PxSiC-camera-video-segment-1743618608072.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSCQSKQN3K%2F20250403%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250403T210046Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIz%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIA0f1nL%2FXh13bkMdIT2AmwaiBdDQqQS9e0S%2FljvQ4gKaAiEA66C2PeoZXPhcqkfZj%2FSmhFW0JZF358cY4pis%2FqXODcgq5gMI9f%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDFOSUKyEHXsANiPAbiq6A14JdRT1c3F1axdrWHYBaP4dVElP6xpcwzyENXLCxholSW0EyRm2KoLDT8a4n9B9ZJuKgjPEt5LjThVMWWKOg40pb0dlf0FjqleeWpY7EqFLWNUnEERekltQN1FrsnhNaybTLBWDUbK7QqIe9ribi4fvhm8GEKsmOeqS1im8xQ8gTbqHK%2BARSuIYRfmjOrPZwuNh9HBM5HM6Sxc9JNYGsc3xbc%2F59xOniBVMkmvIee3HSFyuAGIGoSi4x4FFBPKhucNC%2BHrhhGHM%2F4oXd8tJ2CxfVpjlsGuV3fwyfd452aPIyuOjvS9Ws7cW7VArjg4xRY%2BVM%2BI0%2Fygju8d9lKYak34Yq%2B2TunJCxHrLWr9XmZ1RMg3SsCOdudtodwvZl3mew2CXGL%2B4MOPnu%2B3u1gZD4aKcpbo3ODmhav8qEerJJapJ4ZfbHF%2FH9zwKEGG3Q7ntq1SjpOOItgE%2B6eQ4byWj1pIkxEEUltRHx8Cie8fmzXi0PYYHPkLcAbOGoJPSu6tnhHtgqv7g0yy%2BSo5mqe7G8wZzy3edkd4BpevZaBsLiyQME6wkh%2BxodKB7wbF%2FhoZJ0VY2CJjX1U6i0qYwy8O7vwY6pQH27dqHdKwnvo1S%2B%2BTEkT6eVRCthiib3yqPaq7st7YzJ0LAAIq6zu7k0iBrQLWobAv1vd337f4fGLbgEpafIwWBZouQ11G%2FKzHkRAZN%2F0AhKRdBlD5gDHSwKcdJdJ5aSyQykz0sYCO%2Bq%2BkMt1upPfVIzb48jgPBVMfe7Vbso6c%2FIrl5F9HJojfAmmTbKgieM40Lq0vvS9Shm8v%2FZIOYPJniMHJakrk%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2004%20Apr%202025%2005%3A00%3A46%20GMT&X-Amz-Signature=3f178a38b832aebf63230497c90b7546b32e52e7dd8291947b50b451ee315130'
